$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]92 + "n"
$note = [char]9834

# Column B (numbers) first
$ws.Range("B7").Value = 322
$ws.Range("B8").Value = 325

# Column C (English)
$ws.Range("C7").Value = " Welcome back! " + $note + "[K] Back to your" + $nl + "guild training, eh?"
$ws.Range("C8").Value = " We hope for your continued" + $nl + "patronage of our humble business! " + $note

# Column D (Russian)
$ws.Range("D7").Value = " С возвращением! " + $note + "[K] Снова" + $nl + "тренируетесь в гильдии, да?"
$ws.Range("D8").Value = " Мы надеемся, что вы и дальше" + $nl + "будете поддерживать наше скромное дело! " + $note

# Column E (converted)
$ws.Range("E7").Value = " Ò âïèâñàþåîéåí! " + $note + "[K] Òîïâà" + $nl + "óñåîéñôåóåòû â ãéìûäéé, äà?"
$ws.Range("E8").Value = " Íú îàäååíòÿ, œóï âú é äàìûšå" + $nl + "áôäåóå ðïääåñçéâàóû îàšå òëñïíîïå äåìï! " + $note

$ws.Rows.Item(7).RowHeight = 21.6
$ws.Rows.Item(8).RowHeight = 31.8

$ws.Range("C12").Select() | Out-Null
